$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like labels stored as plain text (shared strings),
# not real dates. Typing "01-07-2021" straight into Value would let Excel's
# smart input turn it into a date serial + a new number-format style, so we
# stage it as a text formula result elsewhere and paste only the value in,
# which preserves the plain-text type without touching any cell's style.
$ws.Range("Z1").Formula = "=""01-07-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A83").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("Z1").Value = ""
$ws.Range("Z1").Clear()

$ws.Range("B83").Value = 765
$ws.Range("C83").Value = 754
$ws.Range("D83").Value = 747
